$wb = $excel.ActiveWorkbook

# Update "Latest Handback DateTime" (column L) for the e58f67ba-... row (row 3)
# on both the zh-cn and de-de localization-status report sheets, reflecting a
# freshly generated handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("L3").Value = "2017-01-03 06:53:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("L3").Value = "2017-01-03 06:53:42"
